$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note in D21 referencing the new shared string
$ws.Range("D21").Value = "*2018-19 same amounts"

# Update the selection to match the post-edit state (C20)
$ws.Range("C20").Select()
